# Corrección de casos de uso y plantilla de casos de uso.
# Se modifican los casos de uso 16 y 17 (CRU egreso y CRU gasto promocional)
# en descripciones, modelo de CU y plantilla de casos de uso; se agrega el
# nuevo caso de uso 25 "Consultar gastos".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# --- CU-17 (fila 21): "CRU egreso" --- (orden elegido para preservar el
# orden de aparición de las nuevas cadenas compartidas, tal como en el commit)
$ws.Range("D21").Value = "CRU egreso"

# --- CU-16 (fila 20): "CRU gasto promocional" ---
$ws.Range("D20").Value = "CRU gasto promocional"
$ws.Range("C20").Value = "El director puede administrar gastos de promociones de facebook"

# --- CU-17 (fila 21) descripción ---
$ws.Range("C21").Value = "El director puede administrar gastos realizados"

# --- Nuevo CU-25 (fila 29): "Consultar gastos" ---
$ws.Range("B29").Value = "CU - 25"
$ws.Range("C29").Value = "El director puede consultar cualquier tipo de gasto de la institución"
$ws.Range("D29").Value = "Consultar gastos"
$ws.Range("E29").Value = "vacio"
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = 1

# --- Actualizar selección de la hoja ---
$ws.Activate()
$ws.Range("D30").Select()
